# Helper: COM shape geometry is exposed in points (1 pt = 12700 EMU) and the
# interop layer keeps that value in a single-precision float, so a plain
# "emu / 12700.0" division can truncate to one EMU below the intended value
# once it is floored back to EMU on save. Nudging by +0.5 EMU before the
# division keeps the float32 round-trip inside the correct integer-EMU
# bucket for the shape sizes used on this deck.
function EmuToPt($emu) {
    return ($emu + 0.5) / 12700
}

$p = $ppt.ActivePresentation

# --- Slide 15: "PM Probe Query for P2MP SR Policy" ---
$s15 = $p.Slides.Item(15)

# Shape 3 "Rectangle 2": widen the code-block textbox and update the figure caption.
$rect2 = $s15.Shapes.Item(3)
$rect2HeightEmu = 2123658
$rect2.TextFrame.AutoSize = 0
$rect2.TextFrame.TextRange.Paragraphs(11).Runs(1).Text = " Figure: Example Probe Query with Replication Segment for P2MP SR Policy"
$rect2.TextFrame.AutoSize = 1
$rect2.Height = EmuToPt $rect2HeightEmu
$rect2.Width = EmuToPt 6819900

# Shape 5 "Rectangle 4": update bullet text; keep the autosized height exactly as-is.
$rect4 = $s15.Shapes.Item(5)
$rect4HeightEmu = 1467005
$rect4.TextFrame.AutoSize = 0
$rect4.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Applicable to one-way mode for delay and loss measurement for P2MP SR Policy."
$rect4.TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "The querier root node sends probe query messages using the Replication Segment for the P2MP SR Policy"
$rect4.TextFrame.TextRange.Paragraphs(4).Runs(1).Text = "This TLV allows the querier root node to identify the responder leaf nodes of the P2MP SR Policy"
$rect4.TextFrame.AutoSize = 1
$rect4.Height = EmuToPt $rect4HeightEmu

# --- Slide 7: "Reverse Path PM" -- TextBox 6 ---
$s7 = $p.Slides.Item(7)
$tb6 = $s7.Shapes.Item(6)
$tb6HeightEmu = 1323439
$tb6.TextFrame.AutoSize = 0
$tb6.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "TLV is mandatory when carried in a probe query message and if responder does not support, it MUST return Error "
$tb6.TextFrame.AutoSize = 1
$tb6.Height = EmuToPt $tb6HeightEmu

# --- Slide 8: "Block Number TLV for Loss Measurement" -- Rectangle 6 ---
$s8 = $p.Slides.Item(8)
$rect6 = $s8.Shapes.Item(5)
$rect6HeightEmu = 1815882
$rect6.TextFrame.AutoSize = 0
$rect6.TextFrame.TextRange.Paragraphs(4).Runs(1).Text = "TLV is mandatory when carried in a probe query message and if responder does not support, it MUST return "
$rect6.TextFrame.AutoSize = 1
$rect6.Height = EmuToPt $rect6HeightEmu
